$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the previous data row (row 4) down into the new
# row 5 so the new row picks up the same cell styles (date format on A,
# boolean style on G, etc.) without introducing any new style entries.
$ws.Range("A4:I4").Copy($ws.Range("A5:I5"))

# Now fill in the new trade's actual values.
$ws.Range("A5").Value = 42636.589039351849
$ws.Range("B5").Value = $false
$ws.Range("C5").Value = 9917.16
$ws.Range("D5").Value = 9948
$ws.Range("E5").Value = 19.29
$ws.Range("F5").Value = 19.41
$ws.Range("G5").Value = $true
$ws.Range("H5").Value = 0.62
$ws.Range("I5").Value = $false
